# Auto-generated Excel COM-interop script
# Applies numeric cell value updates to match the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 396.33334
$ws.Range("I2").Value = 400
$ws.Range("J2").Value = 394.5
$ws.Range("K2").Value = 400
$ws.Range("L2").Value = 394.5
$ws.Range("M2").Value = -287
$ws.Range("N2").Value = -620.5
$ws.Range("H17").Value = 1305752
$ws.Range("J17").Value = 1436177.2
$ws.Range("L17").Value = 4308531.6
$ws.Range("N17").Value = -4308867.6
$ws.Range("H18").Value = 1322.1666
$ws.Range("I18").Value = 310.33334
$ws.Range("J18").Value = 2334
$ws.Range("K18").Value = 310.33334
$ws.Range("L18").Value = 2334
$ws.Range("M18").Value = -26.33334000000002
$ws.Range("N18").Value = -2902
$ws.Range("H33").Value = 106.30769
$ws.Range("I33").Value = 97.2
$ws.Range("K33").Value = 97.2
$ws.Range("M33").Value = 131.8
$ws.Range("H40").Value = 1803.7
$ws.Range("I40").Value = 1638.7333
$ws.Range("J40").Value = 2298.6
$ws.Range("K40").Value = 1638.7333
$ws.Range("L40").Value = 2298.6
$ws.Range("M40").Value = -1463.7333
$ws.Range("N40").Value = -2648.6
$ws.Range("H64").Value = 4187.7085
$ws.Range("I64").Value = 3178.7144
$ws.Range("J64").Value = 5600.3
$ws.Range("K64").Value = 3178.7144
$ws.Range("L64").Value = 5600.3
$ws.Range("M64").Value = -2930.7144
$ws.Range("N64").Value = -6096.3
$ws.Range("H67").Value = 4187.7085
$ws.Range("I67").Value = 3178.7144
$ws.Range("J67").Value = 5600.3
$ws.Range("K67").Value = 3178.7144
$ws.Range("L67").Value = 5600.3
$ws.Range("M67").Value = -2320.7144
$ws.Range("N67").Value = -7316.3
$ws.Range("H112").Value = 17858534
$ws.Range("J112").Value = 20834792
$ws.Range("L112").Value = 62504376
$ws.Range("N112").Value = -62506592
$ws.Range("H113").Value = 2351.4285
$ws.Range("I113").Value = 2351.4285
$ws.Range("K113").Value = 2351.4285
$ws.Range("M113").Value = 902.5715
$ws.Range("H129").Value = 826.3333
$ws.Range("J129").Value = 1427.3
$ws.Range("L129").Value = 4281.9
$ws.Range("N129").Value = -14281.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19244.693
$ws.Range("I32").Value = 2708.1404
$ws.Range("J32").Value = 207761.4
$ws.Range("K32").Value = 2708.1404
$ws.Range("L32").Value = 207761.4
$ws.Range("M32").Value = -2421.1404
$ws.Range("N32").Value = -208335.4
$ws.Range("H61").Value = 2594.457
$ws.Range("I61").Value = 1874.9231
$ws.Range("K61").Value = 1874.9231
$ws.Range("M61").Value = -1662.9231
$ws.Range("H97").Value = 55573304
$ws.Range("I97").Value = 66687664
$ws.Range("J97").Value = 1500
$ws.Range("K97").Value = 66687664
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -66687168
$ws.Range("N97").Value = -2492
$ws.Range("H102").Value = 2100
$ws.Range("I102").Value = 2100
$ws.Range("K102").Value = 2100
$ws.Range("M102").Value = -478
$ws.Range("H112").Value = 12199.667
$ws.Range("J112").Value = 12199.667
$ws.Range("L112").Value = 12199.667
$ws.Range("N112").Value = -15153.667
$ws.Range("H119").Value = 37497
$ws.Range("J119").Value = 37497
$ws.Range("L119").Value = 37497
$ws.Range("N119").Value = -47173
$ws.Range("H122").Value = 1808.0952
$ws.Range("I122").Value = 1765.8462
$ws.Range("K122").Value = 5297.5386
$ws.Range("M122").Value = -2847.5386
$ws.Range("H124").Value = 51485.8
$ws.Range("J124").Value = 51485.8
$ws.Range("L124").Value = 51485.8
$ws.Range("N124").Value = -61305.8
$ws.Range("H125").Value = 28985.715
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 28985.715
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 28985.715
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -38825.715
$ws.Range("H136").Value = 2594.457
$ws.Range("I136").Value = 1874.9231
$ws.Range("K136").Value = 5624.7693
$ws.Range("M136").Value = -3074.7693

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 12245.5
$ws.Range("I86").Value = 2618.5
$ws.Range("J86").Value = 50753.5
$ws.Range("K86").Value = 2618.5
$ws.Range("L86").Value = 50753.5
$ws.Range("M86").Value = -1495.5
$ws.Range("N86").Value = -52999.5
$ws.Range("H89").Value = 12245.5
$ws.Range("I89").Value = 2618.5
$ws.Range("J89").Value = 50753.5
$ws.Range("K89").Value = 13092.5
$ws.Range("L89").Value = 253767.5
$ws.Range("M89").Value = -7476.5
$ws.Range("N89").Value = -264999.5
$ws.Range("H94").Value = 1500.7894
$ws.Range("I94").Value = 1055.6666
$ws.Range("K94").Value = 1055.6666
$ws.Range("M94").Value = -604.6666
$ws.Range("H110").Value = 36666.332
$ws.Range("J110").Value = 36666.332
$ws.Range("L110").Value = 36666.332
$ws.Range("N110").Value = -44846.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4663.75
$ws.Range("I31").Value = 1175.12
$ws.Range("J31").Value = 12592.454
$ws.Range("K31").Value = 1175.12
$ws.Range("L31").Value = 12592.454
$ws.Range("M31").Value = -880.1199999999999
$ws.Range("N31").Value = -13182.454
$ws.Range("H34").Value = 4663.75
$ws.Range("I34").Value = 1175.12
$ws.Range("J34").Value = 12592.454
$ws.Range("K34").Value = 1175.12
$ws.Range("L34").Value = 12592.454
$ws.Range("M34").Value = -973.1199999999999
$ws.Range("N34").Value = -12996.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 546.1667
$ws.Range("I107").Value = 604.2857
$ws.Range("J107").Value = 464.8
$ws.Range("K107").Value = 1812.8571
$ws.Range("L107").Value = 1394.4
$ws.Range("M107").Value = 107.1428999999998
$ws.Range("N107").Value = -5234.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2524.4443
$ws.Range("I80").Value = 2327.8572
$ws.Range("J80").Value = 3212.5
$ws.Range("K80").Value = 2327.8572
$ws.Range("L80").Value = 3212.5
$ws.Range("M80").Value = -1329.8572
$ws.Range("N80").Value = -5208.5
$ws.Range("H83").Value = 2524.4443
$ws.Range("I83").Value = 2327.8572
$ws.Range("J83").Value = 3212.5
$ws.Range("K83").Value = 11639.286
$ws.Range("L83").Value = 16062.5
$ws.Range("M83").Value = -6647.286
$ws.Range("N83").Value = -26046.5
$ws.Range("H132").Value = 2670.2666
$ws.Range("I132").Value = 2376.8
$ws.Range("J132").Value = 3697.4
$ws.Range("K132").Value = 7130.400000000001
$ws.Range("L132").Value = 11092.2
$ws.Range("M132").Value = -4600.400000000001
$ws.Range("N132").Value = -16152.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 862.2
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 903.6667
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 903.6667
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -1493.6667
$ws.Range("H27").Value = 862.2
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 903.6667
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 903.6667
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -1117.6667
$ws.Range("H40").Value = 2930.1177
$ws.Range("I40").Value = 1835.6666
$ws.Range("K40").Value = 1835.6666
$ws.Range("M40").Value = -1699.6666
$ws.Range("H46").Value = 825.2143
$ws.Range("I46").Value = 662.75
$ws.Range("K46").Value = 662.75
$ws.Range("M46").Value = -474.75
$ws.Range("H61").Value = 1174.9166
$ws.Range("I61").Value = 1174.9166
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1174.9166
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -972.9166
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1174.9166
$ws.Range("I113").Value = 1174.9166
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1174.9166
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 995.0834
$ws.Range("N113").ClearContents()
$ws.Range("H127").Value = 49571.668
$ws.Range("J127").Value = 49571.668
$ws.Range("L127").Value = 49571.668
$ws.Range("N127").Value = -59491.668
$ws.Range("H132").Value = 2867.5483
$ws.Range("I132").Value = 1833.3478
$ws.Range("K132").Value = 5500.0434
$ws.Range("M132").Value = -2970.0434

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 161.33333
$ws.Range("I100").Value = 142
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 284
$ws.Range("L100").Value = 400
$ws.Range("M100").Value = 257
$ws.Range("N100").Value = -1482
$ws.Range("H122").Value = 201700.8
$ws.Range("I122").Value = 1000004
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 3000012
$ws.Range("L122").Value = 6375
$ws.Range("M122").Value = -2997562
$ws.Range("N122").Value = -11275
$ws.Range("H123").Value = 32692.154
$ws.Range("J123").Value = 32692.154
$ws.Range("L123").Value = 32692.154
$ws.Range("N123").Value = -42492.15399999999

